$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1931464174454829
$ws.Range("C2").Value = 0.5545171339563862
$ws.Range("J2").Value = 0.009345794392523364
$ws.Range("P2").Value = 0.1401869158878505
$ws.Range("S2").Value = 0.102803738317757
$ws.Range("B3").Value = 0.00558659217877095
$ws.Range("C3").Value = 0.0223463687150838
$ws.Range("J3").Value = 0.0223463687150838
$ws.Range("P3").Value = 0.8044692737430168
$ws.Range("S3").Value = 0.1452513966480447
$ws.Range("J4").Value = 0.04347826086956522
$ws.Range("P4").Value = 0.6956521739130435
$ws.Range("S4").Value = 0.2608695652173913
$ws.Range("B6").Value = 0.07894736842105263
$ws.Range("D6").Value = 0.05263157894736842
$ws.Range("F6").Value = 0.06842105263157895
$ws.Range("J6").Value = 0.1736842105263158
$ws.Range("O6").Value = 0.02105263157894737
$ws.Range("Q6").Value = 0.2
$ws.Range("R6").Value = 0.03684210526315789
$ws.Range("S6").Value = 0.3684210526315789
$ws.Range("B7").Value = 0.1344086021505376
$ws.Range("D7").Value = 0.02688172043010753
$ws.Range("F7").Value = 0.05913978494623656
$ws.Range("J7").Value = 0.1451612903225807
$ws.Range("O7").Value = 0.02150537634408602
$ws.Range("Q7").Value = 0.2150537634408602
$ws.Range("R7").Value = 0.07526881720430108
$ws.Range("S7").Value = 0.3225806451612903
$ws.Range("B8").Value = 0.09116022099447514
$ws.Range("D8").Value = 0.02209944751381215
$ws.Range("F8").Value = 0.06353591160220995
$ws.Range("J8").Value = 0.1104972375690608
$ws.Range("O8").Value = 0.02209944751381215
$ws.Range("Q8").Value = 0.2099447513812155
$ws.Range("R8").Value = 0.07734806629834254
$ws.Range("S8").Value = 0.4033149171270718
$ws.Range("B9").Value = 0.122257053291536
$ws.Range("D9").Value = 0.01253918495297806
$ws.Range("F9").Value = 0.05329153605015674
$ws.Range("J9").Value = 0.109717868338558
$ws.Range("O9").Value = 0.01880877742946709
$ws.Range("Q9").Value = 0.2507836990595611
$ws.Range("R9").Value = 0.05642633228840126
$ws.Range("S9").Value = 0.3761755485893417
$ws.Range("B10").Value = 0.1081081081081081
$ws.Range("D10").Value = 0.01544401544401544
$ws.Range("E10").Value = 0.0007722007722007722
$ws.Range("F10").Value = 0.05714285714285714
$ws.Range("J10").Value = 0.1127413127413127
$ws.Range("O10").Value = 0.02007722007722008
$ws.Range("Q10").Value = 0.2455598455598456
$ws.Range("R10").Value = 0.08030888030888031
$ws.Range("S10").Value = 0.3598455598455598
$ws.Range("G11").Value = 0.1169811320754717
$ws.Range("J11").Value = 0.06792452830188679
$ws.Range("K11").Value = 0.1622641509433962
$ws.Range("L11").Value = 0.6415094339622641
$ws.Range("S11").Value = 0.01132075471698113
$ws.Range("G12").Value = 0.7840909090909091
$ws.Range("J12").Value = 0.1534090909090909
$ws.Range("L12").Value = 0.03409090909090909
$ws.Range("S12").Value = 0.02840909090909091
$ws.Range("G13").Value = 0.6451612903225806
$ws.Range("J13").Value = 0.2903225806451613
$ws.Range("S13").Value = 0.06451612903225806
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.0184331797235023
$ws.Range("H15").Value = 0.1336405529953917
$ws.Range("I15").Value = 0.1059907834101382
$ws.Range("J15").Value = 0.3824884792626728
$ws.Range("K15").Value = 0.06912442396313365
$ws.Range("M15").Value = 0.004608294930875576
$ws.Range("O15").Value = 0.05069124423963134
$ws.Range("S15").Value = 0.2350230414746544
$ws.Range("F16").Value = 0.02752293577981652
$ws.Range("H16").Value = 0.1100917431192661
$ws.Range("I16").Value = 0.1100917431192661
$ws.Range("J16").Value = 0.4541284403669725
$ws.Range("K16").Value = 0.1330275229357798
$ws.Range("M16").Value = 0.009174311926605505
$ws.Range("N16").Value = 0.009174311926605505
$ws.Range("O16").Value = 0.05963302752293578
$ws.Range("S16").Value = 0.0871559633027523
$ws.Range("F17").Value = 0.01094890510948905
$ws.Range("H17").Value = 0.1514598540145985
$ws.Range("I17").Value = 0.1313868613138686
$ws.Range("J17").Value = 0.4306569343065693
$ws.Range("K17").Value = 0.08394160583941605
$ws.Range("M17").Value = 0.02737226277372263
$ws.Range("N17").Value = 0.001824817518248175
$ws.Range("O17").Value = 0.06386861313868614
$ws.Range("S17").Value = 0.09854014598540146
$ws.Range("F18").Value = 0.005882352941176471
$ws.Range("H18").Value = 0.1529411764705882
$ws.Range("I18").Value = 0.1529411764705882
$ws.Range("J18").Value = 0.4764705882352941
$ws.Range("K18").Value = 0.08235294117647059
$ws.Range("M18").Value = 0.01176470588235294
$ws.Range("O18").Value = 0.04117647058823529
$ws.Range("S18").Value = 0.07647058823529412
$ws.Range("F19").Value = 0.01096121416526138
$ws.Range("H19").Value = 0.1711635750421585
$ws.Range("I19").Value = 0.1450252951096122
$ws.Range("J19").Value = 0.3937605396290054
$ws.Range("K19").Value = 0.09527824620573355
$ws.Range("M19").Value = 0.01011804384485666
$ws.Range("O19").Value = 0.06576728499156829
$ws.Range("S19").Value = 0.1079258010118044
